$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NA" text values in C3/D3 are replaced with the number 0.
# (This also drops "NA" from the shared-string table since it becomes unused.)
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Move the active cell / selection (the Name Box used for scrolling to a
# cell) from E10 to D6.
$ws.Range("D6").Select()

# Best-effort: also try to match the window geometry recorded in the
# workbook view (xWindow/yWindow/windowWidth/windowHeight). Some hosts
# don't persist these through the object model, but setting them is
# harmless if unsupported.
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 25600
$win.Height = 14960
